$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.442.26'
$ws.Range("D3").Value = '1.850.10'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6274'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2920'
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07753'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").Value = '1.854.10'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.035'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6819'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001073'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.50'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '29.458.60'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.60'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.444'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.86'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1379'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.372'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.462'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05625'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.127'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7067'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.598'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").Value = '1.225.27'
$ws.Range("E37").Value = '  -1.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01794'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.754'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.447'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").Value = '2.011.53'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.96'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.18'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000120'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.198'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4021'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.030'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1155'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.73%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.676'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.05%  '
